$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New mapping table: column A = new_names, column B = old_names,
# now expanded to include hispanic / not-hispanic variants and a
# renamed geography column (geographic_area_name -> state_county).
$rows = @(
    @("new_names", "old_names"),
    @("id", "id"),
    @("state_county", "geographic_area_name"),
    @("total_nlat", "estimate_total_not_hispanic_or_latino"),
    @("total_nlat_white_alone", "estimate_total_not_hispanic_or_latino_white_alone"),
    @("total_nlat_blackaa_alone", "estimate_total_not_hispanic_or_latino_black_or_african_american_alone"),
    @("total_nlat_native_alone", "estimate_total_not_hispanic_or_latino_american_indian_and_alaska_native_alone"),
    @("total_nlat_asian_alone", "estimate_total_not_hispanic_or_latino_asian_alone"),
    @("total_nlat_island_alone", "estimate_total_not_hispanic_or_latino_native_hawaiian_and_other_pacific_islander_alone"),
    @("total_nlat_other_race", "estimate_total_not_hispanic_or_latino_some_other_race_alone"),
    @("total_nlat_2ormore", "estimate_total_not_hispanic_or_latino_two_or_more_races"),
    @("total_nlat_2ormore_io", "estimate_total_not_hispanic_or_latino_two_or_more_races_two_races_including_some_other_race"),
    @("total_nlat_3ormore_io", "estimate_total_not_hispanic_or_latino_two_or_more_races_two_races_excluding_some_other_race_and_three_or_more_races"),
    @("total_lat", "estimate_total_hispanic_or_latino"),
    @("total_lat_white_alone", "estimate_total_hispanic_or_latino_white_alone"),
    @("total_lat_blackaa_alone", "estimate_total_hispanic_or_latino_black_or_african_american_alone"),
    @("total_lat_native_alone", "estimate_total_hispanic_or_latino_american_indian_and_alaska_native_alone"),
    @("total_lat_asian_alone", "estimate_total_hispanic_or_latino_asian_alone"),
    @("total_lat_island_alone", "estimate_total_hispanic_or_latino_native_hawaiian_and_other_pacific_islander_alone"),
    @("total_lat_other_race", "estimate_total_hispanic_or_latino_some_other_race_alone"),
    @("total_lat_2ormore", "estimate_total_hispanic_or_latino_two_or_more_races"),
    @("total_lat_2ormore_io", "estimate_total_hispanic_or_latino_two_or_more_races_two_races_including_some_other_race"),
    @("total_lat_3ormore_io", "estimate_total_hispanic_or_latino_two_or_more_races_two_races_excluding_some_other_race_and_three_or_more_races")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Widen the two columns to fit the longer labels.
$ws.Columns.Item(1).ColumnWidth = 39.5
$ws.Columns.Item(2).ColumnWidth = 41.666666666666664

# Zoom out a touch and move the selection down onto the new data.
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("A5").Select()
